# =====================================================================================
# Weekly data refresh for "Hortaliza, Vega Monumental Concepcion - Betarraga".
# A new latest observation (rows 140/141, date 2021-09-21) is inserted at the top of the
# date-ordered block (rows 140-195); every existing row below shifts down by two rows
# (one "Primera" + one "Segunda" pair per week), and the two rows that fall off the
# bottom of the original range (old 194/195) reappear as new rows 196/197.
# =====================================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for columns D (Fecha), J (Volumen), M (Precio promedio ponderado) and
# P (Precio $/Kg) for existing rows 140-195 (the only columns that change for these rows).
$rowData = @{
    140 = @{ D = 44460; J = 600; M = 650; P = 130 }
    141 = @{ D = 44460; J = 300; M = 500; P = 100 }
    142 = @{ D = 44427; J = 600; M = 650; P = 130 }
    143 = @{ D = 44427; J = 300; M = 500; P = 100 }
    144 = @{ D = 44343; J = 1000; M = 650; P = 130 }
    145 = @{ D = 44343; J = 500; M = 500; P = 100 }
    146 = @{ D = 44280; J = 800; M = 650; P = 130 }
    147 = @{ D = 44280; J = 400; M = 500; P = 100 }
    148 = @{ D = 44390; J = 600; M = 650; P = 130 }
    149 = @{ D = 44390; J = 300; M = 500; P = 100 }
    150 = @{ D = 44386; J = 800; M = 650; P = 130 }
    151 = @{ D = 44386; J = 400; M = 500; P = 100 }
    152 = @{ D = 44308; J = 600; M = 650; P = 130 }
    153 = @{ D = 44308; J = 300; M = 500; P = 100 }
    154 = @{ D = 44264; J = 600; M = 650; P = 130 }
    155 = @{ D = 44264; J = 300; M = 500; P = 100 }
    156 = @{ D = 44196; J = 800; M = 650; P = 130 }
    157 = @{ D = 44196; J = 400; M = 500; P = 100 }
    158 = @{ D = 44243; J = 800; M = 650; P = 130 }
    159 = @{ D = 44243; J = 400; M = 500; P = 100 }
    160 = @{ D = 44252; J = 800; M = 650; P = 130 }
    161 = @{ D = 44252; J = 400; M = 500; P = 100 }
    162 = @{ D = 44166; J = 600; M = 650; P = 130 }
    163 = @{ D = 44166; J = 300; M = 500; P = 100 }
    164 = @{ D = 44168; J = 600; M = 650; P = 130 }
    165 = @{ D = 44168; J = 300; M = 500; P = 100 }
    166 = @{ D = 44316; J = 1000; M = 650; P = 130 }
    167 = @{ D = 44316; J = 500; M = 500; P = 100 }
    168 = @{ D = 44397; J = 600; M = 650; P = 130 }
    169 = @{ D = 44397; J = 300; M = 500; P = 100 }
    170 = @{ D = 44273; J = 600; M = 650; P = 130 }
    171 = @{ D = 44273; J = 300; M = 500; P = 100 }
    172 = @{ D = 44372; J = 600; M = 650; P = 130 }
    173 = @{ D = 44372; J = 300; M = 500; P = 100 }
    174 = @{ D = 44365; J = 600; M = 650; P = 130 }
    175 = @{ D = 44365; J = 300; M = 500; P = 100 }
    176 = @{ D = 44306; J = 600; M = 650; P = 130 }
    177 = @{ D = 44306; J = 300; M = 500; P = 100 }
    178 = @{ D = 44215; J = 800; M = 650; P = 130 }
    179 = @{ D = 44215; J = 400; M = 500; P = 100 }
    180 = @{ D = 44357; J = 800; M = 650; P = 130 }
    181 = @{ D = 44357; J = 400; M = 500; P = 100 }
    182 = @{ D = 44162; J = 800; M = 650; P = 130 }
    183 = @{ D = 44162; J = 400; M = 500; P = 100 }
    184 = @{ D = 44239; J = 700; M = 643; P = 129 }
    185 = @{ D = 44239; J = 300; M = 500; P = 100 }
    186 = @{ D = 44376; J = 600; M = 650; P = 130 }
    187 = @{ D = 44376; J = 300; M = 500; P = 100 }
    188 = @{ D = 44292; J = 600; M = 650; P = 130 }
    189 = @{ D = 44292; J = 300; M = 500; P = 100 }
    190 = @{ D = 44358; J = 600; M = 650; P = 130 }
    191 = @{ D = 44358; J = 300; M = 500; P = 100 }
    192 = @{ D = 44211; J = 600; M = 650; P = 130 }
    193 = @{ D = 44211; J = 300; M = 500; P = 100 }
    194 = @{ D = 44425; J = 600; M = 650; P = 130 }
    195 = @{ D = 44425; J = 300; M = 500; P = 100 }
}

foreach ($r in $rowData.Keys) {
    $d = $rowData[$r]
    $ws.Cells.Item([int]$r, 4).Value = $d.D
    $ws.Cells.Item([int]$r, 10).Value = $d.J
    $ws.Cells.Item([int]$r, 13).Value = $d.M
    $ws.Cells.Item([int]$r, 16).Value = $d.P
}

# Two brand-new rows appended at the bottom (196-197), duplicating the data that used to
# sit in rows 194-195 before the shift (same market/category/quality template + date).

# Row 196
$ws.Cells.Item(196, 1).Value = 11
$ws.Cells.Item(196, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(196, 3).Value = "Bíobío"
$ws.Cells.Item(196, 4).Value = 44323
$ws.Cells.Item(196, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(196, 5).Value = 8
$ws.Cells.Item(196, 6).Value = 100114014
$ws.Cells.Item(196, 7).Value = "Betarraga"
$ws.Cells.Item(196, 8).Value = "Sin especificar"
$ws.Cells.Item(196, 9).Value = "Primera"
$ws.Cells.Item(196, 10).Value = 600
$ws.Cells.Item(196, 11).Value = 600
$ws.Cells.Item(196, 12).Value = 700
$ws.Cells.Item(196, 13).Value = 650
$ws.Cells.Item(196, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(196, 15).Value = "Región Metropolitana"
$ws.Cells.Item(196, 16).Value = 130
$ws.Cells.Item(196, 17).Value = 5
$ws.Cells.Item(196, 18).Value = "Hortaliza"

# Row 197
$ws.Cells.Item(197, 1).Value = 11
$ws.Cells.Item(197, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(197, 3).Value = "Bíobío"
$ws.Cells.Item(197, 4).Value = 44323
$ws.Cells.Item(197, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(197, 5).Value = 8
$ws.Cells.Item(197, 6).Value = 100114014
$ws.Cells.Item(197, 7).Value = "Betarraga"
$ws.Cells.Item(197, 8).Value = "Sin especificar"
$ws.Cells.Item(197, 9).Value = "Segunda"
$ws.Cells.Item(197, 10).Value = 300
$ws.Cells.Item(197, 11).Value = 500
$ws.Cells.Item(197, 12).Value = 500
$ws.Cells.Item(197, 13).Value = 500
$ws.Cells.Item(197, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(197, 15).Value = "Región Metropolitana"
$ws.Cells.Item(197, 16).Value = 100
$ws.Cells.Item(197, 17).Value = 5
$ws.Cells.Item(197, 18).Value = "Hortaliza"
